# Insert a new weekly price record at row 65 ("Terminal Hortofrutícola Agro
# Chillán" - Poroto verde sheet), pushing all existing rows from 65 downward
# by one row (old row 65 -> new row 66, ..., old row 138 -> new row 139).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 65; Excel shifts rows 65..138 down to 66..139
# and extends the used range accordingly.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new record's data.
$ws.Cells.Item(65, 1).Value = 7
$ws.Cells.Item(65, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(65, 3).Value = "Ñuble"
$ws.Cells.Item(65, 4).Value = 45079
$ws.Cells.Item(65, 5).Value = 16
$ws.Cells.Item(65, 6).Value = 100112031
$ws.Cells.Item(65, 7).Value = "Poroto verde"
$ws.Cells.Item(65, 8).Value = "Magnum"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 30
$ws.Cells.Item(65, 11).Value = 25000
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = 25000
$ws.Cells.Item(65, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(65, 15).Value = "Perú"
$ws.Cells.Item(65, 16).Value = 1000
$ws.Cells.Item(65, 17).Value = 25
$ws.Cells.Item(65, 18).Value = "Hortaliza"
